{"js": "// The author corrected two meeting/time-log entries for \"12/9\" from\n// \"1:..pm\" to \"12:..pm\" (a typo fix: 1pm -> 12pm), by inserting the\n// digit \"2\" right after the \"1\" in \"12/9 at 1...\".\n//\n// We locate each unique timestamp string in the document body and\n// replace it with the corrected text, which preserves the existing\n// run formatting (font/size/etc.) of the text being replaced.\n\nconst body = context.document.body;\n\n// 1) \"12/9 at 1:00pm\" -> \"12/9 at 12:00pm\" (Communication Log table)\nconst firstMatches = body.search(\"12/9 at 1:00pm\", { matchCase: true, matchWholeWord: false });\nfirstMatches.load(\"items\");\nawait context.sync();\n\nif (firstMatches.items.length > 0) {\n  firstMatches.items[0].insertText(\"12/9 at 12:00pm\", \"Replace\");\n  await context.sync();\n}\n\n// 2) \"12/9 at 1:10\" -> \"12/9 at 12:10\" (Personal Time Logs table; the\n// trailing \"pm\" lives in a separate, unchanged run right after this text)\nconst secondMatches = body.search(\"12/9 at 1:10\", { matchCase: true, matchWholeWord: false });\nsecondMatches.load(\"items\");\nawait context.sync();\n\nif (secondMatches.items.length > 0) {\n  secondMatches.items[0].insertText(\"12/9 at 12:10\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# The author corrected two meeting/time-log entries for \"12/9\" from\n# \"1:..pm\" to \"12:..pm\" (a typo fix: 1pm -> 12pm), by inserting the\n# digit \"2\" right after the \"1\" in \"12/9 at 1...\".\n#\n# We use Find/Replace scoped to each unique timestamp string so only\n# that text changes; the surrounding run formatting (font/size/etc.)\n# is preserved because Replace re-uses the formatting of the text it\n# replaces.\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceOne = 1\n\n# 1) \"12/9 at 1:00pm\" -> \"12/9 at 12:00pm\" (Communication Log table)\n$rng1 = $d.Content\n$rng1.Find.Execute(\"12/9 at 1:00pm\", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"12/9 at 12:00pm\", $wdReplaceOne)\n\n# 2) \"12/9 at 1:10\" -> \"12/9 at 12:10\" (Personal Time Logs table; the\n# trailing \"pm\" lives in a separate, unchanged run right after this text)\n$rng2 = $d.Content\n$rng2.Find.Execute(\"12/9 at 1:10\", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"12/9 at 12:10\", $wdReplaceOne)\n"}
